# CRMS-2395 Update Create date column in excel
#
# The "Create Date" column's merge-field placeholder (row 2, column AC)
# is renamed from {vendor:create_date} to {vendor:sf_create_date so the
# generated list pulls the date from the new SalesForce-backed field.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the merge-field token used by the "Create Date" column.
$ws.Range("AC2").Value = "{vendor:sf_create_date"

# Reflect the edit location/selection, matching the user clicking into
# the cell they just changed.
$ws.Range("AC2").Select()
